# Append two new data rows (7 and 8) to the active worksheet, matching
# the records added upstream. Existing data in rows 1-6 is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 7 - Garnlav (Alectoria sarmentosa)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 131255264
$ws.Range("B7").Value = 79244
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("P7").Value = "Stora Vittjärnen, Dlr"
$ws.Range("Q7").Value = 485818
$ws.Range("R7").Value = 6666268
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = "Dalarna"
$ws.Range("U7").Value = "Ludvika"
$ws.Range("V7").Value = "Dalarna"
$ws.Range("W7").Value = "Grangärde"

# Force the date/time-looking values to stay plain text (otherwise Excel
# auto-converts recognizable date/time strings into date serial numbers).
$ws.Range("Y7:AB7").NumberFormat = "@"
$ws.Range("Y7").Value = "2026-02-22"
$ws.Range("Z7").Value = "09:18"
$ws.Range("AA7").Value = "2026-02-22"
$ws.Range("AB7").Value = "09:18"

$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = "Tobias Hellberg"
$ws.Range("AX7").Value = "Tobias Hellberg"

# ---------------------------------------------------------------------
# Row 8 - Tretåig hackspett (Picoides tridactylus)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 131255137
$ws.Range("B8").Value = 57884
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("M8").Value = "äldre spår"
$ws.Range("P8").Value = "Stora Vittjärnen, Dlr"
$ws.Range("Q8").Value = 485867
$ws.Range("R8").Value = 6666265
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = "Dalarna"
$ws.Range("U8").Value = "Ludvika"
$ws.Range("V8").Value = "Dalarna"
$ws.Range("W8").Value = "Grangärde"

$ws.Range("Y8:AB8").NumberFormat = "@"
$ws.Range("Y8").Value = "2026-02-22"
$ws.Range("Z8").Value = "09:12"
$ws.Range("AA8").Value = "2026-02-22"
$ws.Range("AB8").Value = "09:12"

$ws.Range("AC8").Value = "Äldre ringhack på gran"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = "Tobias Hellberg"
$ws.Range("AX8").Value = "Tobias Hellberg"
